# Update 'want to go' counts (column F) across sheets, per commit
# 'Update gh-pages to output generated at 456a3b4'

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 88
$ws.Range("F3").Value = 169
$ws.Range("F5").Value = 3352
$ws.Range("F6").Value = 1168
$ws.Range("F7").Value = 2270
$ws.Range("F12").Value = 1714
$ws.Range("F18").Value = 1631
$ws.Range("F19").Value = 285
$ws.Range("F20").Value = 1345
$ws.Range("F21").Value = 760
$ws.Range("F22").Value = 288
$ws.Range("F23").Value = 640
$ws.Range("F24").Value = 12452
$ws.Range("F25").Value = 12508
$ws.Range("F26").Value = 926
$ws.Range("F29").Value = 272
$ws.Range("F31").Value = 421
$ws.Range("F36").Value = 634

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("F10").Value = 59

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 88
$ws.Range("F4").Value = 169
$ws.Range("F6").Value = 3352
$ws.Range("F7").Value = 1168
$ws.Range("F8").Value = 2270
$ws.Range("F14").Value = 1714
$ws.Range("F23").Value = 1631
$ws.Range("F24").Value = 285
$ws.Range("F25").Value = 1345
$ws.Range("F26").Value = 760
$ws.Range("F27").Value = 288
$ws.Range("F28").Value = 2
$ws.Range("F29").Value = 640
$ws.Range("F30").Value = 12452
$ws.Range("F31").Value = 12508
$ws.Range("F32").Value = 926
$ws.Range("F35").Value = 272
$ws.Range("F36").Value = 38
$ws.Range("F37").Value = 421
$ws.Range("F38").Value = 3
$ws.Range("F46").Value = 634
$ws.Range("F47").Value = 59
